# Advance the cleaning-plan rotation by one week:
#  - "Calendar" sheet: each week's date moves forward 7 days, and the task
#    assignments shift so that week N gets what used to be week N+1's
#    assignments; a new final week is appended with a fresh assignment.
#  - "Roles" sheet: header row content is unchanged (Floor / Kitchen /
#    Bathrooms / Management), only included here for completeness/safety.

$wb = $excel.ActiveWorkbook
$calendar = $wb.Worksheets.Item("Calendar")
$roles = $wb.Worksheets.Item("Roles")

# New state for rows 2..16 of the Calendar sheet:
# (row, date-serial, Arman, Cesare, Claudio, Jaspar, Lea, Mara)
$data = @(
    @(2, 45348, "Management", "Bathrooms", "Floor", "Vacation", "Kitchen", "Vacation"),
    @(3, 45355, "Vacation", "Vacation", "Kitchen", "Management", "Floor", "Bathrooms"),
    @(4, 45362, "Bathrooms", "Management", "Vacation", "Kitchen", "Vacation", "Floor"),
    @(5, 45369, "Kitchen", "Vacation", "Vacation", "Bathrooms", "Vacation", "Management"),
    @(6, 45376, "Anarchy", "Anarchy", "Anarchy", "Anarchy", "Anarchy", "Anarchy"),
    @(7, 45383, "Vacation", "Management", "Kitchen", "Vacation", "Bathrooms", "Floor"),
    @(8, 45390, "Vacation", "Kitchen", "Bathrooms", "Floor", "Management", "Vacation"),
    @(9, 45397, "Management", "Bathrooms", "Floor", "Vacation", "Kitchen", "Vacation"),
    @(10, 45404, "Bathrooms", "Floor", "Vacation", "Management", "Vacation", "Kitchen"),
    @(11, 45411, "Management", "Vacation", "Vacation", "Kitchen", "Floor", "Bathrooms"),
    @(12, 45418, "Floor", "Vacation", "Management", "Bathrooms", "Vacation", "Kitchen"),
    @(13, 45425, "Floor", "Kitchen", "Management", "Vacation", "Bathrooms", "Vacation"),
    @(14, 45432, "Vacation", "Bathrooms", "Kitchen", "Floor", "Management", "Vacation"),
    @(15, 45439, "Vacation", "Floor", "Bathrooms", "Vacation", "Kitchen", "Management"),
    @(16, 45446, "Kitchen", "Vacation", "Floor", "Management", "Vacation", "Bathrooms")
)

foreach ($row in $data) {
    $r = $row[0]
    $calendar.Range("A" + $r).Value = $row[1]
    $calendar.Range("B" + $r).Value = $row[2]
    $calendar.Range("C" + $r).Value = $row[3]
    $calendar.Range("D" + $r).Value = $row[4]
    $calendar.Range("E" + $r).Value = $row[5]
    $calendar.Range("F" + $r).Value = $row[6]
    $calendar.Range("G" + $r).Value = $row[7]
}

# Roles sheet header row stays semantically the same (Floor / Kitchen /
# Bathrooms / Management) - rewrite explicitly so the textual content is
# guaranteed to match even though the diff only reflects a shared-string
# table reshuffle.
$roles.Range("A1").Value = "Floor"
$roles.Range("B1").Value = "Kitchen"
$roles.Range("C1").Value = "Bathrooms"
$roles.Range("D1").Value = "Management"
